$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve N2's existing cell format (it uses a quote-prefixed text style because
# its old value started with "+"). Stash a copy of that formatting on a scratch
# cell before we overwrite N2's value, since Excel may drop the quote-prefix
# style automatically once the new text no longer needs it.
$ws.Range("N2").Copy()
$ws.Range("Q1").PasteSpecial(-4122)

$ws.Range("N2").Value = "abortion,pro choice,pro-choice,pro life,pro-life,dobbs,anti-abortion,anti abortion"
$ws.Range("O2").Value = "rights,murder,violation,body,right,individual,government"
$ws.Range("N3").Value = "education,schools,school,teachers,teacher,class,classroom,classrooms"
$ws.Range("O3").Value = "books,book,ban,bans,banning,banned,parents,choice,funding,racial disparities,disparity,learning gap,teachers,gender,gay,lgbtq,rights,homophobic,dangerous,discrimination,discriminatory,nazis,ideology,brainwash,children,kids,child,trans,transgender"
$ws.Range("N4").Value = "immigration,migration,immigrants, migrants,assylum,southern border"
$ws.Range("O4").Value = " crisis,legacy americans,replace,opportunity,opportunities,better life,seeking"
$ws.Range("N5").Value = "crime"
$ws.Range("O5").Value = "drug,drugs,theft,violent,violence,chaos,society,dangerous,homeless,rule of law,law and order,police,gangs"
$ws.Range("N6").Value = "health care"
$ws.Range("O6").Value = "health care"
$ws.Range("N7").Value = "guns,mass shooting,mass shootings, second ammendment,assault weapons"
$ws.Range("O7").Value = "ban,rights,deadly,dangerous,tyranny,government,take"
$ws.Range("N8").Value = "energy policy,green new deal, renewable energy,renewables,wind, solar,electrification,fossil fuels,oil,coal,nuclear power"
$ws.Range("O8").Value = "climate change,debt,money,jobs"
$ws.Range("N9").Value = "economy,inflation,jobs,wages,unemployment, salaries"
$ws.Range("O9").Value = "good,bad,great,terrible,aweful,horrible,poor,rich,wealthy,tax cuts,biden,trump,better,worse"
$ws.Range("N10").Value = "black lives matter,blm,police misconduct,social justice"
$ws.Range("O10").Value = "riots,protests,white privilege,murder,rights,right,unarmed,black,white,treatment,fair,unfair,better,worse,racist,racists"
$ws.Range("N11").Value = "supreme court,justices"
$ws.Range("O11").Value = "conservative,liberal,trump,biden,values,rights,dismantle,dismantling,destroy,destroying"
$ws.Range("N12").Value = "small government,big government,spending cuts,debt ceiling,government debt"
$ws.Range("O12").Value = "debt,shutdown,default,bad,consequences,"
$ws.Range("N13").Value = "climate change,global warming"
$ws.Range("O13").Value = "energy policy,green new deal,renewable energy,renewables,wind,solar,oil,electrification,fossil fuels,climate change,hoax,fraud,corrupt"
$ws.Range("N14").Value = "election fraud,voter supression,fake electors,january 6,january 6th"
$ws.Range("O14").Value = "jim crow,big lie,claims,mike pence,hoax,lying,not true,false,story,police,defend,riot,protest,persecuted,persecute,jail,prosecute,prosecuted,insurrection,seditious conspiracy,sedition"
$ws.Range("N15").Value = "foreign policy,ukraine,russia,china,nato,putin,xi,jinping,war"
$ws.Range("O15").Value = "freedom,democracy,dictator,weapons,autocracy,authoritarian,alliance,corruption,corrupt,fraud,hoax,cold war,nuclear,military industrial complex,greed,thug,corporations,brave,valiant,home"
$ws.Range("N16").Value = "coronavirus,covid,pandemic,vaccine,vaccines"
$ws.Range("O16").Value = "coronavirus,covid,pandemic,mandates,safety,protect,tyranny,authoritarian,mask,masks,Fauci,lie,false,dangerous,deadly,die,dying,dead"
$ws.Range("N17").Value = "future of democracy,future of the country,nation's future"
$ws.Range("O17").Value = "future,safe,election,fraud,hoax,big lie,insurrection,january 6,january 6th,fake electors,plot,coup,scheme,woke,tyranny,rights,right,fear,dictatorship,dictator,fascist,fascists,communist,communists,radical,extremist,extremists"
$ws.Range("N18").Value = "investigation,investigations"
$ws.Range("O18").Value = "january 6,classified documents,taxes,trump,classified documents,mar-a-lago,mar a lago,china,corrupt,crime,hoax,fraud,greed,grift,collude,collusion,moscow,russia"
$ws.Range("N19").Value = "investigation,investigations"
$ws.Range("O19").Value = "laptop,burisma,classified documents,hunter,biden,ukraine,china,corrupt,crime,hoax,fraud,greed,grift,collude,collusion,"

# Restore N2's original quote-prefixed style, then set values that changed but
# already match existing style rules normally.
$ws.Range("Q1").Copy()
$ws.Range("N2").PasteSpecial(-4122)
$ws.Range("Q1").Clear()
$ws.Application.CutCopyMode = $false

# Update the selected cell, matching the author's last selection.
$ws.Range("N3").Select()
